$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (new content sourced from old row 4) ---
$ws.Range("A3").Value = 131047025
$ws.Range("B3").Value = 89194
$ws.Range("E3").Value = 510
$ws.Range("F3").Value = "Doftskinn"
$ws.Range("G3").Value = "Cystostereum murrayi"
$ws.Range("H3").Value = "(Berk. & M.A.Curtis.) Pouzar"
$ws.Range("Q3").Value = 402314
$ws.Range("R3").Value = 6818423
$ws.Range("Z3").Value = "16:05"
$ws.Range("AB3").Value = "16:05"
# --- Row 4 (new content sourced from old row 3) ---
$ws.Range("A4").Value = 131046847
$ws.Range("B4").Value = 79244
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 402380
$ws.Range("R4").Value = 6818405
$ws.Range("Z4").Value = "17:01"
$ws.Range("AB4").Value = "17:01"
# --- Row 9 (new content sourced from old row 10) ---
$ws.Range("A9").Value = 131046844
$ws.Range("B9").Value = 79244
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 402484
$ws.Range("R9").Value = 6818538
$ws.Range("Z9").Value = "16:23"
$ws.Range("AB9").Value = "16:23"
# --- Row 10 (new content sourced from old row 9) ---
$ws.Range("A10").Value = 131046711
$ws.Range("B10").Value = 83224
$ws.Range("E10").Value = 6440
$ws.Range("F10").Value = "Vitgrynig nållav"
$ws.Range("G10").Value = "Chaenotheca subroscida"
$ws.Range("H10").Value = "(Eitner) Zahlbr."
$ws.Range("Q10").Value = 402363
$ws.Range("R10").Value = 6818428
$ws.Range("Z10").Value = "16:09"
$ws.Range("AB10").Value = "16:09"
# --- Row 11 (new content sourced from old row 12) ---
$ws.Range("A11").Value = 131046735
$ws.Range("B11").Value = 57884
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("M11").Value = "nyligen använt bo"
$ws.Range("Q11").Value = 402448
$ws.Range("R11").Value = 6818295
$ws.Range("Z11").Value = "16:54"
$ws.Range("AB11").Value = "16:54"
# --- Row 12 (new content sourced from old row 13) ---
$ws.Range("A12").Value = 131046788
$ws.Range("M12").Value = "färska spår"
$ws.Range("Q12").Value = 402473
$ws.Range("R12").Value = 6818425
$ws.Range("Z12").Value = "16:47"
$ws.Range("AB12").Value = "16:47"
$ws.Range("AC12").Value = "Färska ringhack (gran)"
# --- Row 13 (new content sourced from old row 11) ---
$ws.Range("A13").Value = 131046763
$ws.Range("B13").Value = 92268
$ws.Range("D13").Value = "VU"
$ws.Range("E13").Value = 1209
$ws.Range("F13").Value = "Rynkskinn"
$ws.Range("G13").Value = "Hermanssonia centrifuga"
$ws.Range("H13").Value = "(P. Karst.) Zmitr."
$ws.Range("M13").ClearContents()
$ws.Range("Q13").Value = 402378
$ws.Range("R13").Value = 6818392
$ws.Range("Z13").Value = "17:02"
$ws.Range("AB13").Value = "17:02"
$ws.Range("AC13").ClearContents()
# --- Row 23 (new content sourced from old row 24) ---
$ws.Range("A23").Value = 131047016
$ws.Range("B23").Value = 57884
$ws.Range("E23").Value = 100109
$ws.Range("F23").Value = "Tretåig hackspett"
$ws.Range("G23").Value = "Picoides tridactylus"
$ws.Range("H23").Value = "(Linnaeus, 1758)"
$ws.Range("M23").Value = "färska spår"
$ws.Range("Q23").Value = 402474
$ws.Range("R23").Value = 6818507
$ws.Range("Z23").Value = "16:22"
$ws.Range("AB23").Value = "16:22"
$ws.Range("AC23").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AE23").Value = $true
# --- Row 24 (new content sourced from old row 23) ---
$ws.Range("A24").Value = 131046845
$ws.Range("B24").Value = 79244
$ws.Range("E24").Value = 6425
$ws.Range("F24").Value = "Garnlav"
$ws.Range("G24").Value = "Alectoria sarmentosa"
$ws.Range("H24").Value = "(Ach.) Ach."
$ws.Range("M24").ClearContents()
$ws.Range("Q24").Value = 402575
$ws.Range("R24").Value = 6818545
$ws.Range("Z24").Value = "16:34"
$ws.Range("AB24").Value = "16:34"
$ws.Range("AC24").ClearContents()
$ws.Range("AE24").Value = $false
